# Lernraten_EEM_Paper.xlsx - "Add files via upload" edit
#
# Summary of the change (from the OOXML diff):
#  - Sheet "Ratio_Total_remanufacturing": column A (Total Additions in MW)
#    is rewritten from hard literal values into formulas that subtract
#    260000 from (a truncated copy of) the previous literal value. This
#    ripples into column C (Ratio in %), which already held a formula.
#  - Sheet "Capacity_LR": column C holds literal (not formula) copies of
#    the same ratios, so those literals are refreshed to the newly
#    recomputed numbers.
#  - Selection / active-sheet bookkeeping changes: the user ends up with
#    the Capacity_LR sheet active or selected (cell G12), having last
#    selected C2:C15 on Ratio_Total_remanufacturing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Ratio_Total_remanufacturing: rewrite column A as formulas (=<val>-260000)
# ---------------------------------------------------------------------
$wsRatio = $wb.Worksheets.Item("Ratio_Total_remanufacturing")

# A2 gets its own (non-shared) formula
$wsRatio.Range("A2").Formula = "=2009687.11397106-260000"

# A3:A6 were filled together -> becomes one shared-formula group
$wsRatio.Range("A3:A6").Formula = "=2009687.11397106-260000"

# A7:A15 each carry their own distinct literal-derived formula
$wsRatio.Range("A7").Formula = "=2009629.96937205-260000"
$wsRatio.Range("A8").Formula = "=2009627.56327918-260000"
$wsRatio.Range("A9").Formula = "=2009614.86082608-260000"
$wsRatio.Range("A10").Formula = "=2009612.04561704-260000"
$wsRatio.Range("A11").Formula = "=2009627.71023508-260000"
$wsRatio.Range("A12").Formula = "=2005593.1895522-260000"
$wsRatio.Range("A13").Formula = "=2006496.92829087-260000"
$wsRatio.Range("A14").Formula = "=2074520.08799963-260000"
$wsRatio.Range("A15").Formula = "=2106807-260000"

# Column C (ratio %) recalculates automatically from the new column A,
# since C already holds "=(B/A)*100"-style formulas.

# ---------------------------------------------------------------------
# 2. Capacity_LR: refresh the literal copies of the ratio column (C)
# ---------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capacity_LR")

$wsCap.Range("C7").Value = 13.148754589667186
$wsCap.Range("C8").Value = 13.48271016363492
$wsCap.Range("C9").Value = 15.019590527250449
$wsCap.Range("C10").Value = 15.331464416467178
$wsCap.Range("C11").Value = 19.45601422569278
$wsCap.Range("C12").Value = 31.339646561083285
$wsCap.Range("C13").Value = 36.637348204567402
$wsCap.Range("C14").Value = 65.024068281366993
$wsCap.Range("C15").Value = 65.629630383683832

# ---------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping
# ---------------------------------------------------------------------
# User's last action on Ratio_Total_remanufacturing: select C2:C15.
$wsRatio.Range("C2:C15").Select()

# Workbook ends with Capacity_LR as the active sheet, cell G12 selected.
$wsCap.Select()
$wsCap.Range("G12").Select()
